$d = $word.ActiveDocument

# --- Locate the last paragraph in the body (ends with "...makes the game more unique. ") ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tail = $lastPara.Range
$tail.Collapse(0)

# --- New blank paragraph (matches the blank spacer paragraphs used throughout the doc) ---
$tail.InsertParagraphAfter()
$pBlank = $d.Paragraphs($d.Paragraphs.Count)

# --- "Planetary Interactions:" heading-style paragraph ---
$pBlank.Range.Collapse(0)
$pBlank.Range.InsertParagraphAfter()
$pHeading = $d.Paragraphs($d.Paragraphs.Count)
$pHeading.Range.Text = "Planetary Interactions:"

# --- First bullet item: sets up the new bulleted list (numId 2) ---
$pHeading.Range.Collapse(0)
$pHeading.Range.InsertParagraphAfter()
$pBullet1 = $d.Paragraphs($d.Paragraphs.Count)
$pBullet1.Range.Text = "Planet interactions have to work in order for players to gather resources otherwise there wouldn’t be much to on them other than explore them in Pulsar. So, there for to create the planet interactions the land alteration factors must be implemented to allow the players to mine, excavate and cut down trees. Also, plant trees and place building parts to make it work."

$gallery = $word.ListGalleries.Item(1)
$bulletTemplate = $gallery.ListTemplates.Item(1)
$pBullet1.Range.ListFormat.ApplyListTemplate($bulletTemplate)

# --- Second bullet item (continues the same list, numId 2) ---
$endOfBullet1 = $d.Paragraphs($d.Paragraphs.Count).Range
$endOfBullet1.Collapse(0)
$endOfBullet1.InsertParagraphAfter()
$pBullet2 = $d.Paragraphs($d.Paragraphs.Count)
$pBullet2.Range.Text = "The terrain can be manipulated by several different factors such as, Humans, Animals & Natural Disasters. Humans could dig down into and generated cave system to gather resources."

# Split the sentence into its two original runs by toggling bold off/on over the
# second sentence only (no visible formatting change, but forces a run boundary).
$searchRange = $pBullet2.Range.Duplicate
$found = $searchRange.Find.Execute("could dig down into and generated cave system to gather resources.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $searchRange.Bold = 1
    $searchRange.Bold = 0
}

# --- Third bullet item (continues the same list, numId 2) ---
$endOfBullet2 = $d.Paragraphs($d.Paragraphs.Count).Range
$endOfBullet2.Collapse(0)
$endOfBullet2.InsertParagraphAfter()
$pBullet3 = $d.Paragraphs($d.Paragraphs.Count)
$pBullet3.Range.Text = "When the players arrive at their destination e.g. a planet, there could be the possibility of a cut scene where there is teleportation systems to the planet if it is well developed and has a good amount of infrastructure the planet itself."
